$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 443.33334
$ws.Range("I61").Value = 432
$ws.Range("K61").Value = 1296
$ws.Range("M61").Value = -1124
$ws.Range("H62").Value = 68763.5
$ws.Range("I62").Value = 201779.6
$ws.Range("J62").Value = 8301.637000000001
$ws.Range("K62").Value = 201779.6
$ws.Range("L62").Value = 8301.637000000001
$ws.Range("M62").Value = -201155.6
$ws.Range("N62").Value = -9549.637000000001
$ws.Range("H65").Value = 68763.5
$ws.Range("I65").Value = 201779.6
$ws.Range("J65").Value = 8301.637000000001
$ws.Range("K65").Value = 1008898
$ws.Range("L65").Value = 41508.185
$ws.Range("M65").Value = -1005778
$ws.Range("N65").Value = -47748.185
$ws.Range("H80").Value = 422.26666
$ws.Range("I80").Value = 237.5
$ws.Range("J80").Value = 791.8
$ws.Range("K80").Value = 712.5
$ws.Range("L80").Value = 2375.4
$ws.Range("M80").Value = 285.5
$ws.Range("N80").Value = -4371.4
$ws.Range("H83").Value = 422.26666
$ws.Range("I83").Value = 237.5
$ws.Range("J83").Value = 791.8
$ws.Range("K83").Value = 2137.5
$ws.Range("L83").Value = 7126.2
$ws.Range("M83").Value = 2854.5
$ws.Range("N83").Value = -17110.2
$ws.Range("H141").Value = 14728.272
$ws.Range("I141").Value = 6201.1
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 18603.3
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = -13423.3
$ws.Range("N141").Value = -310360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7503.794
$ws.Range("I32").Value = 6581.25
$ws.Range("K32").Value = 6581.25
$ws.Range("M32").Value = -6294.25
$ws.Range("H45").Value = 5330587
$ws.Range("I45").Value = 7194018
$ws.Range("J45").Value = 6498.7144
$ws.Range("K45").Value = 7194018
$ws.Range("L45").Value = 6498.7144
$ws.Range("M45").Value = -7193641
$ws.Range("N45").Value = -7252.7144
$ws.Range("H74").Value = 49685.53
$ws.Range("I74").Value = 4693.931
$ws.Range("K74").Value = 4693.931
$ws.Range("M74").Value = -3819.931
$ws.Range("H77").Value = 49685.53
$ws.Range("I77").Value = 4693.931
$ws.Range("K77").Value = 23469.655
$ws.Range("M77").Value = -19101.655
$ws.Range("H122").Value = 1306088.5
$ws.Range("I122").Value = 3802.4285
$ws.Range("J122").Value = 2318977.5
$ws.Range("K122").Value = 11407.2855
$ws.Range("L122").Value = 6956932.5
$ws.Range("M122").Value = -8957.2855
$ws.Range("N122").Value = -6961832.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30306190
$ws.Range("J20").Value = 9997.5
$ws.Range("L20").Value = 9997.5
$ws.Range("N20").Value = -10491.5
$ws.Range("H134").Value = 2744.9788
$ws.Range("I134").Value = 1430.5385
$ws.Range("J134").Value = 9152.875
$ws.Range("K134").Value = 4291.6155
$ws.Range("L134").Value = 27458.625
$ws.Range("M134").Value = -1756.6155
$ws.Range("N134").Value = -32528.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3480.65
$ws.Range("I31").Value = 1542.9445
$ws.Range("J31").Value = 3906
$ws.Range("K31").Value = 1542.9445
$ws.Range("L31").Value = 3906
$ws.Range("M31").Value = -1247.9445
$ws.Range("N31").Value = -4496
$ws.Range("H34").Value = 3480.65
$ws.Range("I34").Value = 1542.9445
$ws.Range("J34").Value = 3906
$ws.Range("K34").Value = 1542.9445
$ws.Range("L34").Value = 3906
$ws.Range("M34").Value = -1340.9445
$ws.Range("N34").Value = -4310
$ws.Range("H58").Value = 2164.3845
$ws.Range("I58").Value = 2061.375
$ws.Range("K58").Value = 2061.375
$ws.Range("M58").Value = -1858.375
$ws.Range("H122").Value = 3519.889
$ws.Range("I122").Value = 3542.5
$ws.Range("J122").Value = 3501.8
$ws.Range("K122").Value = 10627.5
$ws.Range("L122").Value = 10505.4
$ws.Range("M122").Value = -8177.5
$ws.Range("N122").Value = -15405.4
$ws.Range("H132").Value = 93572
$ws.Range("I132").Value = 113265.78
$ws.Range("K132").Value = 339797.34
$ws.Range("M132").Value = -337267.34
$ws.Range("H136").Value = 2164.3845
$ws.Range("I136").Value = 2061.375
$ws.Range("K136").Value = 6184.125
$ws.Range("M136").Value = -3634.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1879.6305
$ws.Range("I113").Value = 3077
$ws.Range("K113").Value = 9231
$ws.Range("M113").Value = -7061
$ws.Range("H132").Value = 1308.4762
$ws.Range("J132").Value = 1382.4445
$ws.Range("L132").Value = 12442.0005
$ws.Range("N132").Value = -17502.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3317494.2
$ws.Range("I102").Value = 4446067.5
$ws.Range("K102").Value = 4446067.5
$ws.Range("M102").Value = -4444445.5
$ws.Range("H122").Value = 743000.5
$ws.Range("I122").Value = 989610.1
$ws.Range("J122").Value = 3171.6667
$ws.Range("K122").Value = 2968830.3
$ws.Range("L122").Value = 9515.000100000001
$ws.Range("M122").Value = -2966380.3
$ws.Range("N122").Value = -14415.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4141.95
$ws.Range("I7").Value = 2878.5
$ws.Range("J7").Value = 9195.75
$ws.Range("K7").Value = 2878.5
$ws.Range("L7").Value = 9195.75
$ws.Range("M7").Value = -2766.5
$ws.Range("N7").Value = -9419.75
$ws.Range("H16").Value = 1704.2307
$ws.Range("I16").Value = 1587.125
$ws.Range("J16").Value = 1891.6
$ws.Range("K16").Value = 1587.125
$ws.Range("L16").Value = 1891.6
$ws.Range("M16").Value = -1417.125
$ws.Range("N16").Value = -2231.6
$ws.Range("H68").Value = 2409.6
$ws.Range("I68").Value = 1955.75
$ws.Range("J68").Value = 4225
$ws.Range("K68").Value = 1955.75
$ws.Range("L68").Value = 4225
$ws.Range("M68").Value = -1206.75
$ws.Range("N68").Value = -5723
$ws.Range("H71").Value = 2409.6
$ws.Range("I71").Value = 1955.75
$ws.Range("J71").Value = 4225
$ws.Range("K71").Value = 9778.75
$ws.Range("L71").Value = 21125
$ws.Range("M71").Value = -6034.75
$ws.Range("N71").Value = -28613
$ws.Range("H122").Value = 5618.1904
$ws.Range("I122").Value = 4063.6667
$ws.Range("J122").Value = 7690.8887
$ws.Range("K122").Value = 12191.0001
$ws.Range("L122").Value = 23072.6661
$ws.Range("M122").Value = -9741.000100000001
$ws.Range("N122").Value = -27972.6661
$ws.Range("H126").Value = 4141.95
$ws.Range("I126").Value = 2878.5
$ws.Range("J126").Value = 9195.75
$ws.Range("K126").Value = 8635.5
$ws.Range("L126").Value = 27587.25
$ws.Range("M126").Value = -6165.5
$ws.Range("N126").Value = -32527.25
$ws.Range("H132").Value = 3005.05
$ws.Range("I132").Value = 2227.8333
$ws.Range("K132").Value = 6683.499899999999
$ws.Range("M132").Value = -4153.499899999999
$ws.Range("H136").Value = 32680.53
$ws.Range("I136").Value = 41074.73
$ws.Range("K136").Value = 123224.19
$ws.Range("M136").Value = -120674.19

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4304.4375
$ws.Range("I122").Value = 4315.1
$ws.Range("K122").Value = 12945.3
$ws.Range("M122").Value = -10495.3
